# Natmi following Dr Hou advice
# Update LR-pairs data: recompute row 2, replace row 3 contents (FAPs target),
# and append two new target-cluster rows (M1, sCs) as rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pomc"
$ws.Cells.Item(2, 3).Value = "Mc5r"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.55477
$ws.Cells.Item(2, 8).Value = 4.66431
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.4635846666666667
$ws.Cells.Item(2, 14).Value = 1.390754
$ws.Cells.Item(2, 15).Value = 0.1767302775232392
$ws.Cells.Item(2, 16).Value = 0.1862343266337822
$ws.Cells.Item(2, 17).Value = 0.7207675321933332
$ws.Cells.Item(2, 18).Value = 6.486907789739999
$ws.Cells.Item(2, 19).Value = 0.1767302775232392
$ws.Cells.Item(2, 20).Value = 0.1862343266337822

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pomc"
$ws.Cells.Item(3, 3).Value = "Mc5r"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.55477
$ws.Cells.Item(3, 8).Value = 4.66431
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.735013
$ws.Cells.Item(3, 14).Value = 5.205038999999999
$ws.Cells.Item(3, 15).Value = 0.6614311279991165
$ws.Cells.Item(3, 16).Value = 0.6970010032454158
$ws.Cells.Item(3, 17).Value = 2.697546162009999
$ws.Cells.Item(3, 18).Value = 24.27791545808999
$ws.Cells.Item(3, 19).Value = 0.6614311279991165
$ws.Cells.Item(3, 20).Value = 0.6970010032454158

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pomc"
$ws.Cells.Item(4, 3).Value = "Mc5r"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.55477
$ws.Cells.Item(4, 8).Value = 4.66431
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.022927
$ws.Cells.Item(4, 14).Value = 0.06878099999999999
$ws.Cells.Item(4, 15).Value = 0.008740356107784638
$ws.Cells.Item(4, 16).Value = 0.00921038747341239
$ws.Cells.Item(4, 17).Value = 0.03564621178999999
$ws.Cells.Item(4, 18).Value = 0.3208159061099999
$ws.Cells.Item(4, 19).Value = 0.008740356107784638
$ws.Cells.Item(4, 20).Value = 0.00921038747341239

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Pomc"
$ws.Cells.Item(5, 3).Value = "Mc5r"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.55477
$ws.Cells.Item(5, 8).Value = 4.66431
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.401595
$ws.Cells.Item(5, 14).Value = 0.80319
$ws.Cells.Item(5, 15).Value = 0.1530982383698596
$ws.Cells.Item(5, 16).Value = 0.1075542826473895
$ws.Cells.Item(5, 17).Value = 0.6243878581499999
$ws.Cells.Item(5, 18).Value = 3.746327148899999
$ws.Cells.Item(5, 19).Value = 0.1530982383698596
$ws.Cells.Item(5, 20).Value = 0.1075542826473895
